$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.878.47'
$ws.Range("E2").Value = '  -2.26%  '
$ws.Range("D3").Value = '1.656.41'
$ws.Range("E3").Value = '  -0.67%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.19'
$ws.Range("E5").Value = '  -0.60%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.00%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3895'
$ws.Range("E7").Value = '  -1.69%  '
$ws.Range("E8").Value = '  -3.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '51.56'
$ws.Range("E9").Value = '  -0.94%  '
$ws.Range("E10").Value = '  -2.99%  '
$ws.Range("E11").Value = '  +0.07%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08491'
$ws.Range("E12").Value = '  -1.03%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '24.05'
$ws.Range("E13").Value = '  -2.08%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.068'
$ws.Range("E14").Value = '  -3.49%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.134'
$ws.Range("E15").Value = '  +2.37%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001315'
$ws.Range("E16").Value = '  -1.80%  '
$ws.Range("D17").Value = '1.654.69'
$ws.Range("E17").Value = '  -0.22%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '94.26'
$ws.Range("E18").Value = '  -1.00%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07014'
$ws.Range("E19").Value = '  +0.11%  '
$ws.Range("E20").Value = '  -4.49%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.992'
$ws.Range("E21").Value = '  -0.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  +0.27%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.77'
$ws.Range("E23").Value = '  +0.07%  '
$ws.Range("D24").Value = '23.868.63'
$ws.Range("E24").Value = '  -2.27%  '
$ws.Range("E25").Value = '  -0.81%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.967'
$ws.Range("E26").Value = '  -3.98%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.11'
$ws.Range("E27").Value = '  -1.99%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '154.00'
$ws.Range("E28").Value = '  -2.42%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.439'
$ws.Range("E29").Value = '  -0.08%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '138.25'
$ws.Range("E30").Value = '  -3.21%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.877'
$ws.Range("E31").Value = '  -2.50%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.495'
$ws.Range("E32").Value = '  -2.04%  '
$ws.Range("D33").Value = '1.835.46'
$ws.Range("E33").Value = '  -0.29%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.016'
$ws.Range("E34").Value = '  -4.39%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.08203'
$ws.Range("E35").Value = '  -0.71%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02915'
$ws.Range("E36").Value = '  -5.01%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.667'
$ws.Range("E37").Value = '  -3.41%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '10.86'
$ws.Range("E38").Value = '  -2.59%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2681'
$ws.Range("E39").Value = '  -3.20%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.09168'
$ws.Range("E40").Value = '  -1.13%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.7588'
$ws.Range("E41").Value = '  -1.75%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '13.58'
$ws.Range("E42").Value = '  -1.73%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.421'
$ws.Range("E43").Value = '  -1.91%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.61'
$ws.Range("E44").Value = '  +0.68%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6952'
$ws.Range("E45").Value = '  -2.55%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.455'
$ws.Range("E46").Value = '  -3.49%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.104'
$ws.Range("E47").Value = '  -0.85%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.9997'
$ws.Range("E48").Value = '  -0.03%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.08310'
$ws.Range("E49").Value = '  -1.50%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '134.27'
$ws.Range("E50").Value = '  -1.84%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.227'
$ws.Range("E51").Value = '  -3.29%  '
